# REFRESH INFO EXCEL sabado 05-10
#
# The missing round scores for the pending contestants are entered and the
# standings table (A3:L8) is re-ordered by total score (column D) from
# highest to lowest, same as the worksheet's own sort button would do.
# Column D keeps its SUM formula and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SOLISTA NOVEL FEMENINO")
$ws.Activate()

function Set-Row($row, $name, $club, $scores) {
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $club
    $ws.Cells.Item($row, 6).Value  = $scores[0]
    $ws.Cells.Item($row, 7).Value  = $scores[1]
    $ws.Cells.Item($row, 8).Value  = $scores[2]
    $ws.Cells.Item($row, 9).Value  = $scores[3]
    $ws.Cells.Item($row, 10).Value = $scores[4]
    $ws.Cells.Item($row, 11).Value = $scores[5]
    $ws.Cells.Item($row, 12).Value = $scores[6]
}

# Final standings, already sorted by total score (column D) descending,
# exactly as the worksheet ends up after judging is completed for everyone.
Set-Row 3 "SAMANTHA GONZALEZ" "HABANA CARACAS" @(10,7,8,8,7,7,8)
Set-Row 4 "ARIAGNA CARMONA"   "HABANA CARACAS" @(10,7,7,7,6,7,8)
Set-Row 5 "NICOLE ESCALONA"   "INDEPENDIENTE"  @(10,7,7,7,7,7,7)
Set-Row 6 "JENNIFER FARFAN"   "DOMUS DANCE"    @(10,7,7,6,8,6,6)
Set-Row 7 "KAREN MELLADO"     "PASION LATINA"  @(10,7,6,6,6,6,7)
Set-Row 8 "SARA CASTRO"       "PASION LATINA"  @(9,6,6,6,6,7,7)

# --- Re-apply the highlight formatting that ends up on L4 after the refresh ---
$srcSheet = $wb.Worksheets.Item("SOLISTA NOVEL MASCULINO")
$srcSheet.Range("H5").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Update the last active selection on the sheet ---
$ws.Range("F8").Select() | Out-Null
